$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 47.17718133333333
$ws.Range("H2").Value = 141.531544
$ws.Range("I2").Value = 0.9278465172287805
$ws.Range("J2").Value = 0.9278465172287805
$ws.Range("M2").Value = 0.4890553333333333
$ws.Range("N2").Value = 1.467166
$ws.Range("O2").Value = 0.9644476581758422
$ws.Range("P2").Value = 0.9644476581758422
$ws.Range("Q2").Value = 23.07225214270044
$ws.Range("R2").Value = 207.650269284304
$ws.Range("S2").Value = 0.8948594006879086
$ws.Range("T2").Value = 0.8948594006879086
$ws.Range("G3").Value = 47.17718133333333
$ws.Range("H3").Value = 141.531544
$ws.Range("I3").Value = 0.9278465172287805
$ws.Range("J3").Value = 0.9278465172287805
$ws.Range("O3").Value = 0.03555234182415776
$ws.Range("P3").Value = 0.03555234182415776
$ws.Range("Q3").Value = 0.8505102250773332
$ws.Range("R3").Value = 7.654592025696
$ws.Range("S3").Value = 0.03298711654087189
$ws.Range("T3").Value = 0.03298711654087189
$ws.Range("I4").Value = 0.02782443874781146
$ws.Range("J4").Value = 0.02782443874781146
$ws.Range("M4").Value = 0.4890553333333333
$ws.Range("N4").Value = 1.467166
$ws.Range("O4").Value = 0.9644476581758422
$ws.Range("P4").Value = 0.9644476581758422
$ws.Range("Q4").Value = 0.6918951082944444
$ws.Range("R4").Value = 6.22705597465
$ws.Range("S4").Value = 0.02683521479038393
$ws.Range("T4").Value = 0.02683521479038392
$ws.Range("I5").Value = 0.02782443874781146
$ws.Range("J5").Value = 0.02782443874781146
$ws.Range("O5").Value = 0.03555234182415776
$ws.Range("P5").Value = 0.03555234182415776
$ws.Range("S5").Value = 0.0009892239574275333
$ws.Range("T5").Value = 0.0009892239574275331
$ws.Range("H6").Value = 6.761849000000001
$ws.Range("I6").Value = 0.04432904402340805
$ws.Range("J6").Value = 0.04432904402340804
$ws.Range("M6").Value = 0.4890553333333333
$ws.Range("N6").Value = 1.467166
$ws.Range("O6").Value = 0.9644476581758422
$ws.Range("P6").Value = 0.9644476581758422
$ws.Range("Q6").Value = 1.102306105548222
$ws.Range("R6").Value = 9.920754949934
$ws.Range("S6").Value = 0.04275304269754971
$ws.Range("T6").Value = 0.0427530426975497
$ws.Range("H7").Value = 6.761849000000001
$ws.Range("I7").Value = 0.04432904402340805
$ws.Range("J7").Value = 0.04432904402340804
$ws.Range("O7").Value = 0.03555234182415776
$ws.Range("P7").Value = 0.03555234182415776
$ws.Range("S7").Value = 0.001576001325858341
$ws.Range("T7").Value = 0.001576001325858341
